$d = $word.ActiveDocument

# Locate the target paragraph via its distinctive (soon-to-be-old) opening text.
$rng = $d.Content
$rng.Find.ClearFormatting()
$found = $rng.Find.Execute("Kubalulekile enganeni yakho ukuba nabangani.")
if (-not $found) {
    throw "Could not locate target paragraph"
}

$para = $rng.Paragraphs(1).Range

# Rebuild the paragraph with the two new runs (translated Zulu text), replacing
# the old 3-run layout (text / shaded space / text) with a clean 2-run layout
# and no shading, per the target edit.
$frag = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="00000074"><w:pPr><w:spacing w:after="240" w:before="240" w:line="240" w:lineRule="auto"/><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r><w:t xml:space="preserve">Kubalulekile enganeni yakho ukuba nabangani. Lokhu kubasiza ekutheni  babe abamngani abahle bakwazi nokuqonda imizwa yabanye abantu. Kuphinde kube ingxenye yokukhula.</w:t></w:r><w:r><w:t xml:space="preserve">Ingane yakho ingase ifune ukuzizwasengathi ihlangene nontanga yayo, futhi yenze nezinto abanye abacabanga ukuthi zinhle. Bangakhuluma noma baqgoke njenga banye ngesikhsthi besa zama ukuzithola ukuthi bangobani. Sometimes, when teens have a hard time making friends, they might be picked on or bullied by other kids. </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$para.InsertXML($frag)
